$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap L1/M1 headers ---
$ws.Range("L1").Value = "Encargado"
$ws.Range("M1").Value = "Ruta Imagen"

# --- Fix existing rows 2 & 3 ---
$ws.Range("J2").Value = 25568.75052516204
$ws.Range("B3").Value = 123456
$ws.Range("J3").Value = 25568.75052516204

# --- Row 4 (new) ---
$ws.Range("A4").Value = "'"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 14019
$ws.Range("C4").Value = "Diego"
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "Calderon"
$ws.Range("F4").Value = "Davila"
$ws.Range("G4").Value = "Jefe de Departamento"
$ws.Range("H4").Value = "Jefe de la Unidad de Informatica"
$ws.Range("I4").Value = "23/03/2023"
$ws.Range("J4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J4").Value = 25568.7505250463
$ws.Range("K4").Value = 14253
$ws.Range("L4").Value = "Niels"
$ws.Range("M4").Value = "C:/Users/MrJua/Downloads/104115574_3109942045730233_5820694040960332009_n.jpg"

# --- Row 5 (new) ---
$ws.Range("A5").Value = "'"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 156465
$ws.Range("C5").Value = "Carlos"
$ws.Range("D5").Value = "Alberto"
$ws.Range("E5").Value = "Catarino"
$ws.Range("F5").Value = "Corralco"
$ws.Range("G5").Value = "Sub-director"
$ws.Range("H5").Value = "Subdirector Academico"
$ws.Range("I5").Value = "23/03/2023"
$ws.Range("J5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J5").Value = 25568.75052494213
$ws.Range("K5").Value = 123123
$ws.Range("L5").Value = "Niels"
$ws.Range("M5").Value = "C:/Users/MrJua/Pictures/Haruu.png"

# --- Row 6 (new) ---
$ws.Range("A6").Value = "'"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 123456
$ws.Range("C6").Value = "Miguel"
$ws.Range("D6").Value = "Angel"
$ws.Range("E6").Value = "Elizondo"
$ws.Range("F6").Value = "Herrera"
$ws.Range("G6").Value = "Sub-director"
$ws.Range("H6").Value = "Subdirector Academico"
$ws.Range("I6").Value = "23/03/2023"
$ws.Range("J6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J6").Value = 45374
$ws.Range("K6").Value = 789541
$ws.Range("L6").Value = "Niels"
$ws.Range("M6").Value = "C:/Users/MrJua/Desktop/SCI/Fotos/XXMKYX_00.jpeg"

# --- Row 7 (new) ---
$ws.Range("A7").Value = "'"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 456987
$ws.Range("C7").Value = "Jose"
$ws.Range("D7").Value = "Angel"
$ws.Range("E7").Value = "De la Rosa"
$ws.Range("F7").Value = "Aviles"
$ws.Range("G7").Value = "Director"
$ws.Range("H7").Value = "Director de la Unidad Academica"
$ws.Range("I7").Value = "23/03/2023"
$ws.Range("J7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J7").Value = 45374
$ws.Range("K7").Value = 142537
$ws.Range("L7").Value = "Mike"
$ws.Range("M7").Value = "C:/Users/MrJua/Downloads/Bocchi/83kcuvo4ad851.png"

# --- Row 8 (new) ---
$ws.Range("A8").Value = "'"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 15014
$ws.Range("C8").Value = "Roberto"
$ws.Range("D8").Value = "Ivan"
$ws.Range("E8").Value = "De la Rosa"
$ws.Range("F8").Value = "Aviles"
$ws.Range("G8").Value = "Director"
$ws.Range("H8").Value = "Director de la Unidad Academica"
$ws.Range("I8").Value = "23/03/2023"
$ws.Range("J8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J8").Value = 45091
$ws.Range("K8").Value = 1548947
$ws.Range("L8").Value = "Juan"
$ws.Range("M8").Value = "C:/Users/MrJua/Downloads/Bocchi/OneDrive-2022-05-04/378.png"

# --- Row 9 (new) ---
$ws.Range("A9").Value = "'"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 14018
$ws.Range("C9").Value = "Juan"
$ws.Range("D9").Value = "Carlos"
$ws.Range("E9").Value = "Calderon"
$ws.Range("F9").Value = "Davila"
$ws.Range("G9").Value = "Director"
$ws.Range("H9").Value = "Director de la Unidad Academica"
$ws.Range("I9").Value = "23/03/2023"
$ws.Range("J9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J9").Value = 45374
$ws.Range("K9").Value = 48456621
$ws.Range("L9").Value = "Niels"
$ws.Range("M9").Value = "C:/Users/MrJua/Downloads/yo.jpg"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "'"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 140198
$ws.Range("C10").Value = "Juan"
$ws.Range("D10").Value = "Carlos"
$ws.Range("E10").Value = "Calderon"
$ws.Range("F10").Value = "Davila"
$ws.Range("G10").Value = "Jefe de Departamento"
$ws.Range("H10").Value = "Jefe de la Unidad de Informatica"
$ws.Range("I10").Value = "23/03/2023"
$ws.Range("J10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J10").Value = 45374
$ws.Range("K10").Value = 14125574
$ws.Range("L10").Value = "Niels"
$ws.Range("M10").Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"

# --- Row 11 (new) ---
$ws.Range("A11").Value = "'"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 123456
$ws.Range("C11").Value = "Juan"
$ws.Range("D11").Value = "Carlos"
$ws.Range("E11").Value = "Calderon"
$ws.Range("F11").Value = "Davila"
$ws.Range("G11").Value = "Director"
$ws.Range("H11").Value = "Director de la Unidad Academica"
$ws.Range("I11").Value = "23/03/2023"
$ws.Range("J11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J11").Value = 45374
$ws.Range("K11").Value = 498748
$ws.Range("L11").Value = "Niels"
$ws.Range("M11").Value = "C:/Users/MrJua/Downloads/Imagen de WhatsApp 2022-12-13 a las 03.15.27.jpg"

# --- Row 12 (new) ---
$ws.Range("A12").Value = "'"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "'56649874"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "Jose"
$ws.Range("D12").Value = "Angel"
$ws.Range("E12").Value = "Hernandez"
$ws.Range("F12").Value = "Olguin"
$ws.Range("G12").Value = "Director"
$ws.Range("H12").Value = "Director de la Unidad Academica"
$ws.Range("I12").Value = "23/03/2023"
$ws.Range("J12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J12").Value = 45374
$ws.Range("K12").Value = 140198
$ws.Range("L12").Value = "Juan Carlos Calderon"
$ws.Range("M12").Value = "C:/Users/MrJua/Downloads/Imagen de WhatsApp 2022-12-13 a las 03.15.27.jpg"
